# Weekly update: insert the latest week's price record for
# Agrícola del Norte S.A. de Arica - Zapallo (Camote) ahead of the
# existing history, pushing the prior rows (old 7..18) down to (8..19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7; Excel shifts rows 7..18 down to 8..19
# and copies the formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with this week's record.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44483
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112045
$ws.Range("G7").Value = "Zapallo"
$ws.Range("H7").Value = "Camote"
$ws.Range("I7").Value = "1a nueva(o)"
$ws.Range("J7").Value = 1300
$ws.Range("K7").Value = 550
$ws.Range("L7").Value = 580
$ws.Range("M7").Value = 565
$ws.Range("N7").Value = "$/kilo (volumen en unidades)"
$ws.Range("O7").Value = "Perú"
$ws.Range("P7").Value = 565
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
